# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "FRANCISCO LUIS CASTILLO JACINTO" (doc 9185297) records move up to
# rows 16-19 (periods sorted ascending 2003..2006), and the
# "LEWIS NADITH MENDOZA BOSSA" (doc 7938359) record moves down to row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block for B16:G20 (Tipo Doc is "CC" for every row, unchanged)
$data = @(
    @{ row = 16; doc = "9185297";  nombre = "FRANCISCO LUIS CASTILLO JACINTO"; periodo = "2003"; mora = 35112; salario = 877803 },
    @{ row = 17; doc = "9185297";  nombre = "FRANCISCO LUIS CASTILLO JACINTO"; periodo = "2004"; mora = 35112; salario = 877803 },
    @{ row = 18; doc = "9185297";  nombre = "FRANCISCO LUIS CASTILLO JACINTO"; periodo = "2005"; mora = 35112; salario = 877803 },
    @{ row = 19; doc = "9185297";  nombre = "FRANCISCO LUIS CASTILLO JACINTO"; periodo = "2006"; mora = 35112; salario = 877803 },
    @{ row = 20; doc = "7938359";  nombre = "LEWIS NADITH MENDOZA BOSSA";      periodo = "2412"; mora = 15600; salario = 1300000 }
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Range("C$r").Value = $item.doc
    $ws.Range("D$r").Value = $item.nombre
    $ws.Range("E$r").Value = $item.periodo
    $ws.Range("F$r").Value = $item.mora
    $ws.Range("G$r").Value = $item.salario
}
